$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Report handback for xinjiang/7ed42eb6-27ad-4a83-baec-771b15d51314.md
# and xinjiang/ffff5977e3da-bc1e-43b4-b9a8-4e91e21a3592.md
# on both the zh-cn and de-de locale sheets: status moves from
# "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Target File" / "Latest Handback File" columns get populated,
# and "Latest Handback DateTime" is stamped.
# ----------------------------------------------------------------------

$hyperlinkColor = 15570276  # COM BGR encoding of RGB(0x64,0x95,0xED) == OOXML rgb FF6495ED

function Set-HandbackRow($ws, [int]$row, [string]$srcDisplay, [string]$srcUrl, [string]$xlfDisplay, [string]$xlfUrl, [string]$handbackDateTime) {
    # Status: handed back, in sync with source
    $ws.Cells.Item($row, 2).Value = "Handed back: in sync with en-US"

    # Latest Target File (column E) - points at the source markdown file
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $srcDisplay
    $ws.Hyperlinks.Add($eCell, $srcUrl, "", "", $srcDisplay)
    $eCell.Font.Underline = $true
    $eCell.Font.Color = $hyperlinkColor

    # Latest Handback File (column F) - points at the handed-back xlf
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $xlfDisplay
    $ws.Hyperlinks.Add($fCell, $xlfUrl, "", "", $xlfDisplay)
    $fCell.Font.Underline = $true
    $fCell.Font.Color = $hyperlinkColor

    # Latest Handback DateTime (column G)
    $ws.Cells.Item($row, 7).Value = $handbackDateTime
}

# ---------------------- zh-cn sheet ----------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/b296278569f11e568e596ea73f827c5a29dec3ce/e2e/7ed42eb6-27ad-4a83-baec-771b15d51314.md"
$zhXlfDisplay = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eeebc1ca5f55a587ee570d2dc7238283f0e0709f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.zh-cn.xlf"

Set-HandbackRow $wsZh 2 "7ed42eb6-27ad-4a83-baec-771b15d51314.md" $mdUrl $zhXlfDisplay $zhXlfUrl "2016-02-17 03:22:25"
Set-HandbackRow $wsZh 3 "7ed42eb6-27ad-4a83-baec-771b15d51314.md" $mdUrl $zhXlfDisplay $zhXlfUrl "2016-02-17 03:22:25"

# ---------------------- de-de sheet ----------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlfDisplay = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46e1d5934b6692a728db1ccfc0fccd845bed2853/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.de-de.xlf"

Set-HandbackRow $wsDe 2 "7ed42eb6-27ad-4a83-baec-771b15d51314.md" $mdUrl $deXlfDisplay $deXlfUrl "2016-02-17 03:22:42"
Set-HandbackRow $wsDe 3 "7ed42eb6-27ad-4a83-baec-771b15d51314.md" $mdUrl $deXlfDisplay $deXlfUrl "2016-02-17 03:22:42"
